$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like "1.00" or
# "0.0000258" or "95.604.75" are not reinterpreted as numbers/dates by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '95.604.75'
$ws.Range("E2").Value = '  -0.62%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.475.87'
$ws.Range("E3").Value = '  +4.56%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.49'
$ws.Range("E5").Value = '  -2.92%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '646.64'
$ws.Range("E6").Value = '  -0.72%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.47'
$ws.Range("E7").Value = '  +7.91%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.404'
$ws.Range("E8").Value = '  -2.99%  '

# Row 9
$ws.Range("E9").Value = '  +0.05%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("E10").Value = '  +1.78%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.475.95'
$ws.Range("E11").Value = '  +4.65%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.198'
$ws.Range("E12").Value = '  -3.29%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.92'
$ws.Range("E13").Value = '  +4.81%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.11'
$ws.Range("E14").Value = '  +0.29%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '95.210.45'
$ws.Range("E15").Value = '  -0.72%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.134.46'
$ws.Range("E16").Value = '  +4.83%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000258'
$ws.Range("E17").Value = '  +3.20%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.51'
$ws.Range("E18").Value = '  +0.64%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.451.60'
$ws.Range("E19").Value = '  +3.79%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.05'
$ws.Range("E20").Value = '  +6.87%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.69'
$ws.Range("E21").Value = '  +12.10%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.516'
$ws.Range("E22").Value = '  +4.84%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '505.57'
$ws.Range("E23").Value = '  +0.74%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.20'
$ws.Range("E24").Value = '  -3.49%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000193'
$ws.Range("E25").Value = '  -1.90%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.65'
$ws.Range("E26").Value = '  +3.06%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.02'
$ws.Range("E27").Value = '  +0.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.13'
$ws.Range("E28").Value = '  +1.49%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.660.50'
$ws.Range("E29").Value = '  +4.83%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.75'
$ws.Range("E30").Value = '  +7.12%  '

# Row 31
$ws.Range("E31").Value = '  +0.00%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.76'
$ws.Range("E32").Value = '  +12.70%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.138'
$ws.Range("E33").Value = '  -2.93%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.184'
$ws.Range("E34").Value = '  -0.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '31.21'
$ws.Range("E35").Value = '  +12.20%  '

# Row 36
$ws.Range("E36").Value = '  +0.17%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.573'
$ws.Range("E37").Value = '  +6.16%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.83'
$ws.Range("E38").Value = '  +4.57%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.45'
$ws.Range("E39").Value = '  -1.70%  '

# Row 40
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '522.66'
$ws.Range("E40").Value = '  +3.80%  '

# Row 41
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.01%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.151'
$ws.Range("E42").Value = '  +0.97%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.918'
$ws.Range("E43").Value = '  +11.29%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.13'
$ws.Range("E44").Value = '  -0.73%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.71'
$ws.Range("E45").Value = '  +4.50%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0418'
$ws.Range("E46").Value = '  +2.71%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.62'
$ws.Range("E47").Value = '  +2.76%  '

# Row 48
$ws.Range("E48").Value = '  -3.51%  '

# Row 49
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.16'
$ws.Range("E49").Value = '  +10.56%  '

# Row 50
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.64'
$ws.Range("E50").Value = '  +1.50%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.21'
$ws.Range("E51").Value = '  +3.06%  '
